# Update the Project Planner / Gantt chart sheet:
#  - rename the project title
#  - replace the generic "Activity NN" placeholder rows with the real
#    Nutritional Food Database Project work-breakdown-structure rows
#  - clear out the now-unused trailing rows (27-39)
#  - hide the Actual Start / Actual Duration / Percent Complete columns
#    (E:G) since this plan only tracks the "Plan" columns for now
#  - refresh the sheet's zoom level and selected cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# --- Project title -----------------------------------------------------
$ws.Range("B1").Value = "Nutritional Food Database Project"

# --- Activity table (Activity name, Plan Start, Plan Duration) ---------
$activities = @(
    @("1.1 Assign Project manager",        1,  1),
    @("1.2 Develop Project Charter",        2,  2),
    @("1.3 Meet Stake Holders",             4,  3),
    @("2.1 Collect Requirement",            7,  2),
    @("2.2 Establish Scope Plan",           9,  3),
    @("2.3 Create WBS",                    12,  1),
    @("2.4 Create Procurement Plan",       13,  2),
    @("2.5 Create Gantt Chart",            15,  1),
    @("2.6 Estimate Activity Resources & Cost", 16, 2),
    @("3.1 Develop Front-End interface",   18,  3),
    @("3.2 Develop Back-End interface",    21,  3),
    @("3.3 Integrate features",            24,  2),
    @("3.4 Quality Assurance & Testing",   26,  3),
    @("3.5 Deployment",                    29,  3),
    @("4.1 Monitor Project Work",          32,  2),
    @("4.2 Control Scope",                 34,  2),
    @("4.3 Control Schedule",              36,  1),
    @("4.4 Control Resources & Cost",      37,  1),
    @("4.5 Perform Quality Control",       38,  1),
    @("4.6 Report Performance",            39,  2),
    @("5.1 Final Performance Review",      40,  1),
    @("5.2 Prepare Final Report",          41,  2)
)

$startRow = 5
for ($i = 0; $i -lt $activities.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $activities[$i][0]
    $ws.Cells.Item($row, 3).Value = $activities[$i][1]
    $ws.Cells.Item($row, 4).Value = $activities[$i][2]
}

# --- Clear the leftover placeholder rows --------------------------------
$ws.Range("B27:G39").ClearContents()

# --- Hide the Actual/Percent-complete columns, no longer tracked -------
$ws.Range("E1:G1").EntireColumn.Hidden = $true

# --- View tweaks (zoom + selection) -------------------------------------
$excel.ActiveWindow.Zoom = 90
$ws.Range("B28").Select()
